$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 165130
$ws.Range("C4").Value = 156076
$ws.Range("C5").Value = 9054
$ws.Range("C8").Value = 64.95
